$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
# B2: "4" -> "3" (keep as text, matching original inline-string type)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "3"

$ws.Range("D2").Value = 0.10175
$ws.Range("E2").Value = -0.121

$ws.Range("K2").Value = 11.485
$ws.Range("L2").Value = 0.2431717128943469
$ws.Range("M2").Value = 0.5468
$ws.Range("N2").Value = 0.0005481704260651628
$ws.Range("O2").Value = 0.04760992599042228
$ws.Range("P2").Value = 0.5468
$ws.Range("Q2").Value = 0.0005481704260651628
$ws.Range("R2").Value = 0.04760992599042228

$ws.Range("U2").Value = 152.88
$ws.Range("V2").Value = 0.1532631578947369
$ws.Range("W2").Value = 0.01701754385964912
$ws.Range("X2").Value = 0.05301032849526706
$ws.Range("Y2").Value = -0.03599278463561793
$ws.Range("Z2").Value = 0.1026761451335898

$ws.Range("AB2").Value = 0.05121153129360456
$ws.Range("AC2").Value = -0.05121153129360456
$ws.Range("AD2").Value = 388.5

$ws.Range("AF2").Value = 388.5
$ws.Range("AG2").Value = 235.62
$ws.Range("AH2").Value = 0.2803030303030303
$ws.Range("AI2").Value = 0.4072284357606314
$ws.Range("AJ2").Value = 0.1910762942779292
$ws.Range("AK2").Value = 0.2941095702320473

# --- Row 3 ---
# D3, E3 are new cells
$ws.Range("D3").Value = -0.0475
$ws.Range("E3").Value = -0.121

$ws.Range("K3").Value = 26.1
$ws.Range("L3").Value = 0.467741935483871
$ws.Range("M3").Value = 0.001
$ws.Range("N3").Value = 0.000001064056182166418
$ws.Range("O3").Value = 0.00003831417624521073
$ws.Range("P3").Value = 0.001
$ws.Range("Q3").Value = 0.000001064056182166418
$ws.Range("R3").Value = 0.00003831417624521073

$ws.Range("U3").Value = 136.2
$ws.Range("V3").Value = 0.1449244520110662
$ws.Range("W3").Value = 0.1165178571428571
$ws.Range("X3").Value = 0.05255099385128935
$ws.Range("Y3").Value = 0.0639668632915678
$ws.Range("Z3").Value = 0.1458823529411764

$ws.Range("AB3").Value = 0.05012067139234247
$ws.Range("AC3").Value = -0.05012067139234247
$ws.Range("AD3").Value = 330.4

$ws.Range("AF3").Value = 330.4
$ws.Range("AG3").Value = 194.2
$ws.Range("AH3").Value = 0.2601165170839238
$ws.Range("AI3").Value = 0.3894848520570553
$ws.Range("AJ3").Value = 0.1712522045855379
$ws.Range("AK3").Value = 0.272714506389552

# --- Row 4 ---
$ws.Range("D4").Value = 0.251
$ws.Range("E4").ClearContents()

$ws.Range("K4").Value = 0.485
$ws.Range("L4").Value = 0.1600660066006601
$ws.Range("M4").Value = 0.5458
$ws.Range("N4").Value = 0.01344334975369458
$ws.Range("O4").Value = 1.125360824742268
$ws.Range("P4").Value = 0.5458
$ws.Range("Q4").Value = 0.01344334975369458
$ws.Range("R4").Value = 1.125360824742268

$ws.Range("T4").Value = 0

$ws.Range("U4").Value = 4.88
$ws.Range("V4").Value = 0.1201970443349754
$ws.Range("W4").Value = 0.01701754385964912
$ws.Range("X4").Value = 0.05301032849526706
$ws.Range("Y4").Value = -0.03599278463561793
$ws.Range("Z4").Value = 0.0789268038551706

$ws.Range("AB4").Value = 0.05121153129360456
$ws.Range("AC4").Value = -0.05121153129360456
$ws.Range("AD4").Value = 15

$ws.Range("AF4").Value = 15
$ws.Range("AG4").Value = 10.12
$ws.Range("AH4").Value = 0.2697841726618705
$ws.Range("AI4").Value = 0.2814258911819887
$ws.Range("AJ4").Value = 0.1995268138801262
$ws.Range("AK4").Value = 0.2090045435770343

# --- Row 5 (becomes Union Capital Limited, old LankaBangla row replaced) ---
$ws.Range("B5").Value = "Union Capital Limited (DSE:UNIONCAP)"

$ws.Range("D5").ClearContents()
$ws.Range("E5").ClearContents()

$ws.Range("G5").Value = -0
$ws.Range("H5").Value = -0
$ws.Range("I5").Value = -0
$ws.Range("J5").Value = -0

$ws.Range("K5").Value = -15.1
$ws.Range("L5").Value = 1.301724137931034
$ws.Range("M5").Value = -0
$ws.Range("N5").Value = -0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = -0
$ws.Range("Q5").Value = -0
$ws.Range("R5").Value = 0

$ws.Range("T5").ClearContents()

$ws.Range("U5").Value = 11.8
$ws.Range("V5").Value = 0.6900584795321637
$ws.Range("W5").Value = -0.6138211382113821
$ws.Range("X5").Value = 0.1082263329630232
$ws.Range("Y5").Value = -0.7220474711744053
$ws.Range("Z5").Value = -0.2966751918158568

$ws.Range("AB5").Value = 0.06167651009416354
$ws.Range("AC5").Value = -0.06167651009416354
$ws.Range("AD5").Value = 43.1

$ws.Range("AF5").Value = 43.1
$ws.Range("AG5").Value = 31.3
$ws.Range("AH5").Value = 0.7159468438538206
$ws.Range("AI5").Value = 0.8223621446288876
$ws.Range("AJ5").Value = 0.646694214876033
$ws.Range("AK5").Value = 0.7707461216449151

# --- Row 6: remove entirely (Union Capital data moved into row 5) ---
$ws.Rows("6").Delete()
